# Replace leading space before cardinality fragments (" ..N" / " ..*")
# with an underscore ("_..N" / "_..*") across the relevant sheets.

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "Coverage";          Cell = "B5";   Value = "_..1" },
    @{ Sheet = "Device";            Cell = "B2";   Value = "_..1" },
    @{ Sheet = "DocumentReference"; Cell = "B7";   Value = "_..1" },
    @{ Sheet = "Encounter";         Cell = "B4";   Value = "_..*" },
    @{ Sheet = "Location";          Cell = "B2";   Value = "_..*" },
    @{ Sheet = "Location";          Cell = "C3";   Value = "_..1" },
    @{ Sheet = "Observation";       Cell = "AA43"; Value = "_..0" },
    @{ Sheet = "Observation";       Cell = "O51";  Value = "_..0" },
    @{ Sheet = "Organization";      Cell = "B3";   Value = "_..1" },
    @{ Sheet = "Organization";      Cell = "B4";   Value = "_..1" },
    @{ Sheet = "Organization";      Cell = "B5";   Value = "_..4" },
    @{ Sheet = "Organization";      Cell = "B6";   Value = "_..1" },
    @{ Sheet = "Organization";      Cell = "B7";   Value = "_..1" },
    @{ Sheet = "Organization";      Cell = "B8";   Value = "_..1" },
    @{ Sheet = "Organization";      Cell = "B9";   Value = "_..1" },
    @{ Sheet = "Organization";      Cell = "C10";  Value = "_..1" },
    @{ Sheet = "Organization";      Cell = "C13";  Value = "_..1" },
    @{ Sheet = "Patient";           Cell = "C6";   Value = "_..1" },
    @{ Sheet = "Practitioner";      Cell = "B2";   Value = "_..1" },
    @{ Sheet = "Practitioner";      Cell = "B3";   Value = "_..1" },
    @{ Sheet = "Practitioner";      Cell = "B4";   Value = "_..4" },
    @{ Sheet = "Practitioner";      Cell = "B5";   Value = "_..1" },
    @{ Sheet = "Practitioner";      Cell = "B6";   Value = "_..1" },
    @{ Sheet = "Provenance";        Cell = "B2";   Value = "_..*" }
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $ws.Range($change.Cell).Value = $change.Value
}
